$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.5777512135871
$ws.Range("D2").Value = 3.362394463865506
$ws.Range("E2").Value = 11.04860540404796
$ws.Range("F2").Value = 16.79084270565066
$ws.Range("G2").Value = 17.61973111071123
$ws.Range("H2").Value = 10.4599317781811
$ws.Range("I2").Value = 24.77537602406732
$ws.Range("M2").Value = 22.7778303133824
$ws.Range("O2").Value = 14.4824990068163
$ws.Range("B3").Value = 10.09165265870313
$ws.Range("D3").Value = 3.329972848910835
$ws.Range("E3").Value = 11.19866697865833
$ws.Range("F3").Value = 16.56556633730244
$ws.Range("G3").Value = 17.04796391240582
$ws.Range("H3").Value = 10.46809987669079
$ws.Range("I3").Value = 23.92236474027304
$ws.Range("M3").Value = 21.60504172128407
$ws.Range("O3").Value = 14.40037371734788
$ws.Range("B4").Value = 9.77972136880941
$ws.Range("D4").Value = 3.310152322000479
$ws.Range("E4").Value = 11.30055468126692
$ws.Range("F4").Value = 16.43295368607714
$ws.Range("G4").Value = 16.69663508149155
$ws.Range("H4").Value = 10.47602578787155
$ws.Range("I4").Value = 23.38852412075184
$ws.Range("M4").Value = 20.84900971390675
$ws.Range("O4").Value = 14.35549443064932
$ws.Range("B5").Value = 9.649339327792488
$ws.Range("D5").Value = 3.302103992195173
$ws.Range("E5").Value = 11.34449147859242
$ws.Range("F5").Value = 16.38042326429856
$ws.Range("G5").Value = 16.55369682521483
$ws.Range("H5").Value = 10.4799854820251
$ws.Range("I5").Value = 23.16882602430515
$ws.Range("M5").Value = 20.53212769440946
$ws.Range("O5").Value = 14.33861496519153
$ws.Range("B6").Value = 9.6274959769519
$ws.Range("D6").Value = 3.300769508644771
$ws.Range("E6").Value = 11.35193210978269
$ws.Range("F6").Value = 16.37179395842083
$ws.Range("G6").Value = 16.52998496481043
$ws.Range("H6").Value = 10.48068699660504
$ws.Range("I6").Value = 23.13222708914376
$ws.Range("M6").Value = 20.47898665133786
$ws.Range("O6").Value = 14.33589762702497
$ws.Range("B7").Value = 9.777976050616024
$ws.Range("D7").Value = 3.310043654053339
$ws.Range("E7").Value = 11.30113748884144
$ws.Range("F7").Value = 16.43223903121438
$ws.Range("G7").Value = 16.69470602705082
$ws.Range("H7").Value = 10.47607623777376
$ws.Range("I7").Value = 23.38556938509171
$ws.Range("M7").Value = 20.84477137819979
$ws.Range("O7").Value = 14.35526106529229
$ws.Range("B8").Value = 10.41300428212256
$ws.Range("D8").Value = 3.351201722176252
$ws.Range("E8").Value = 11.09830204666442
$ws.Range("F8").Value = 16.71202892356154
$ws.Range("G8").Value = 17.4228131955012
$ws.Range("H8").Value = 10.46214289282142
$ws.Range("I8").Value = 24.48357104998238
$ws.Range("M8").Value = 22.38104315695076
$ws.Range("O8").Value = 14.45304078272462
$ws.Range("B9").Value = 11.54716175994669
$ws.Range("D9").Value = 3.432292822422975
$ws.Range("E9").Value = 10.77948361310689
$ws.Range("F9").Value = 17.30251008306994
$ws.Range("G9").Value = 18.83637384449622
$ws.Range("H9").Value = 10.4579903951963
$ws.Range("I9").Value = 26.54157115539382
$ws.Range("M9").Value = 25.10001136050784
$ws.Range("O9").Value = 14.68801227292028
$ws.Range("B10").Value = 12.30801962228966
$ws.Range("D10").Value = 3.491693821329387
$ws.Range("E10").Value = 10.59557325284718
$ws.Range("F10").Value = 17.75698726042245
$ws.Range("G10").Value = 19.85128155464474
$ws.Range("H10").Value = 10.46914288363796
$ws.Range("I10").Value = 27.97764363572612
$ws.Range("M10").Value = 26.90970817540043
$ws.Range("O10").Value = 14.88576867391357
$ws.Range("B11").Value = 12.63770905254609
$ws.Range("D11").Value = 3.51859276343485
$ws.Range("E11").Value = 10.52331323955317
$ws.Range("F11").Value = 17.96716867180226
$ws.Range("G11").Value = 20.30511405609574
$ws.Range("H11").Value = 10.47730790584003
$ws.Range("I11").Value = 28.61131732846459
$ws.Range("M11").Value = 27.69088005555768
$ws.Range("O11").Value = 14.98087593346541
$ws.Range("B12").Value = 12.76014500190385
$ws.Range("D12").Value = 3.528753948018684
$ws.Range("E12").Value = 10.49762921783234
$ws.Range("F12").Value = 18.04716496312076
$ws.Range("G12").Value = 20.47563968657562
$ws.Range("H12").Value = 10.48084415879364
$ws.Range("I12").Value = 28.84823259360918
$ws.Range("M12").Value = 27.98056068825404
$ws.Range("O12").Value = 15.01760133069132
$ws.Range("B13").Value = 12.73388409595646
$ws.Range("D13").Value = 3.526566773891551
$ws.Range("E13").Value = 10.50308541552896
$ws.Range("F13").Value = 18.0299197190973
$ws.Range("G13").Value = 20.43897612017011
$ws.Range("H13").Value = 10.48006281052536
$ws.Range("I13").Value = 28.79734728132676
$ws.Range("M13").Value = 27.91844663812816
$ws.Range("O13").Value = 15.00966075967768
$ws.Range("B14").Value = 12.64783048013756
$ws.Range("D14").Value = 3.51942926917379
$ws.Range("E14").Value = 10.52116631785135
$ws.Range("F14").Value = 17.97374240183602
$ws.Range("G14").Value = 20.31917105557725
$ws.Range("H14").Value = 10.47758993346734
$ws.Range("I14").Value = 28.63087066794704
$ws.Range("M14").Value = 27.71483556360074
$ws.Range("O14").Value = 14.98388329459962
$ws.Range("B15").Value = 12.59480499716588
$ws.Range("D15").Value = 3.515053892635086
$ws.Range("E15").Value = 10.53246125856297
$ws.Range("F15").Value = 17.93938230750827
$ws.Range("G15").Value = 20.24560804855941
$ws.Range("H15").Value = 10.47613307500881
$ws.Range("I15").Value = 28.52849643074207
$ws.Range("M15").Value = 27.58931716270106
$ws.Range("O15").Value = 14.96818544599858
$ws.Range("B16").Value = 12.28613875786514
$ws.Range("D16").Value = 3.489932821933923
$ws.Range("E16").Value = 10.6005284471497
$ws.Range("F16").Value = 17.74331293927236
$ws.Range("G16").Value = 19.82144678906545
$ws.Range("H16").Value = 10.46867148308131
$ws.Range("I16").Value = 27.935817298648
$ws.Range("M16").Value = 26.85780272447634
$ws.Range("O16").Value = 14.87965425619749
$ws.Range("B17").Value = 12.09253727102384
$ws.Range("D17").Value = 3.474485230819086
$ws.Range("E17").Value = 10.64523404237038
$ws.Range("F17").Value = 17.62384744958011
$ws.Range("G17").Value = 19.55907596705302
$ws.Range("H17").Value = 10.4648858742495
$ws.Range("I17").Value = 27.56703857858405
$ws.Range("M17").Value = 26.39820260456051
$ws.Range("O17").Value = 14.82664088912529
$ws.Range("B18").Value = 11.97963871097446
$ws.Range("D18").Value = 3.465589042752897
$ws.Range("E18").Value = 10.67201833107248
$ws.Range("F18").Value = 17.55546291179498
$ws.Range("G18").Value = 19.40743864005825
$ws.Range("H18").Value = 10.4629995427043
$ws.Range("I18").Value = 27.35309326737244
$ws.Range("M18").Value = 26.12989674278003
$ws.Range("O18").Value = 14.79663535906128
$ws.Range("B19").Value = 11.9411496050057
$ws.Range("D19").Value = 3.462575253889002
$ws.Range("E19").Value = 10.68126973084142
$ws.Range("F19").Value = 17.53236829400953
$ws.Range("G19").Value = 19.35597819594768
$ws.Range("H19").Value = 10.4624108511124
$ws.Range("I19").Value = 27.28034744021061
$ws.Range("M19").Value = 26.03837637081958
$ws.Range("O19").Value = 14.78656041113835
$ws.Range("B20").Value = 12.11330666413242
$ws.Range("D20").Value = 3.476130858974398
$ws.Range("E20").Value = 10.64036394163073
$ws.Range("F20").Value = 17.63653136876093
$ws.Range("G20").Value = 19.58708271882944
$ws.Range("H20").Value = 10.46525873635407
$ws.Range("I20").Value = 27.60648718073081
$ws.Range("M20").Value = 26.44753791497165
$ws.Range("O20").Value = 14.83223412628382
$ws.Range("B21").Value = 12.67317226500905
$ws.Range("D21").Value = 3.521526458696195
$ws.Range("E21").Value = 10.51580962818977
$ws.Range("F21").Value = 17.99023275115675
$ws.Range("G21").Value = 20.35439832984664
$ws.Range("H21").Value = 10.47830422259939
$ws.Range("I21").Value = 28.67985312858521
$ws.Range("M21").Value = 27.77480806761729
$ws.Range("O21").Value = 14.9914357286352
$ws.Range("B22").Value = 13.02501161652216
$ws.Range("D22").Value = 3.551046962394086
$ws.Range("E22").Value = 10.44421247937255
$ws.Range("F22").Value = 18.22371813701005
$ws.Range("G22").Value = 20.84804702027464
$ws.Range("H22").Value = 10.48941967819591
$ws.Range("I22").Value = 29.36353649753372
$ws.Range("M22").Value = 28.60648652072245
$ws.Range("O22").Value = 15.09960857500772
$ws.Range("B23").Value = 12.83852887867384
$ws.Range("D23").Value = 3.535307209351442
$ws.Range("E23").Value = 10.48151490716155
$ws.Range("F23").Value = 18.09891935736739
$ws.Range("G23").Value = 20.58535569324226
$ws.Range("H23").Value = 10.483250419378
$ws.Range("I23").Value = 29.00033978324685
$ws.Range("M23").Value = 28.16589974808649
$ws.Range("O23").Value = 15.04150750436598
$ws.Range("B24").Value = 12.10392178428187
$ws.Range("D24").Value = 3.475386917191627
$ws.Range("E24").Value = 10.64256234497052
$ws.Range("F24").Value = 17.63079603451993
$ws.Range("G24").Value = 19.57442333276932
$ws.Range("H24").Value = 10.46508926179412
$ws.Range("I24").Value = 27.58865846035571
$ws.Range("M24").Value = 26.4252461094112
$ws.Range("O24").Value = 14.82970394945853
$ws.Range("B25").Value = 11.25277425403624
$ws.Range("D25").Value = 3.410357138661289
$ws.Range("E25").Value = 10.85705554993737
$ws.Range("F25").Value = 17.13878579173324
$ws.Range("G25").Value = 18.45711702608127
$ws.Range("H25").Value = 10.45662335764776
$ws.Range("I25").Value = 25.99699407541738
$ws.Range("M25").Value = 24.39693683269266
$ws.Range("O25").Value = 14.61993558176098
